$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "cambios de mayo de mayo" - update the reporting period (row 8) from
# the 2021 Q4 period to the 2022 Q1 period, and bump the validation /
# update dates accordingly.
$ws.Range("A8").Value = 2022
$ws.Range("B8").Value = 44562
$ws.Range("C8").Value = 44651
$ws.Range("S8").Value = 44659
$ws.Range("T8").Value = 44659

# Reflect the new on-screen scroll position / selection left behind by
# the edit (view had scrolled back to the left edge, selection moved to
# C11).
$ws.Range("C11").Select()
